$d = $word.ActiveDocument

function Replace-ParagraphRuns($paraIndex, $innerXml) {
    $p = $d.Paragraphs.Item($paraIndex)
    $full = $p.Range
    $r = $d.Range($full.Start, $full.End - 1)
    $xml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $innerXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($xml)
}

# 1. "This is to certify that" -> split into two runs with gramStart/gramEnd proofErr around "that"
$p1 = '<w:r w:rsidRPr="008C44EF"><w:rPr><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">This is to certify </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>that</w:t></w:r><w:proofErr w:type="gramEnd"/>'

# 2. "Has successfully completed" -> split into two runs with gramStart/gramEnd proofErr around "completed"
$p2 = '<w:r w:rsidRPr="002434C7"><w:rPr><w:bCs/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">Has successfully </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:bCs/><w:szCs w:val="28"/></w:rPr><w:t>completed</w:t></w:r><w:proofErr w:type="gramEnd"/>'

# 3 & 4. "Location LOCATION Date START_DATE END_DATE" -> wrap LOCATION with spellStart/spellEnd; split START_DATE END_DATE into two runs
$p3 = '<w:r w:rsidRPr="002434C7"><w:rPr><w:b/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Location</w:t></w:r><w:r w:rsidRPr="002434C7"><w:rPr><w:b/><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="006A22F2"><w:rPr><w:b/><w:color w:val="636363"/><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr><w:t>LOCATION</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="006A22F2" w:rsidRPr="002434C7"><w:rPr><w:b/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidRPr="002434C7"><w:rPr><w:b/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Date</w:t></w:r><w:r w:rsidRPr="002434C7"><w:rPr><w:b/><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="006A22F2"><w:rPr><w:b/><w:color w:val="636363"/><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr><w:t xml:space="preserve">START_DATE </w:t></w:r><w:r><w:rPr><w:b/><w:color w:val="636363"/><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr><w:t>END_DATE</w:t></w:r>'

# 5. ", Virginia  2" + "0110" -> split into ", " run + gramStart + "Virginia  2" run + "0110" run + gramEnd
$p5 = '<w:r><w:rPr><w:color w:val="636363"/><w:szCs w:val="32"/></w:rPr><w:t>9817</w:t></w:r><w:r w:rsidR="002C1399" w:rsidRPr="008C44EF"><w:rPr><w:color w:val="636363"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:color w:val="636363"/><w:szCs w:val="32"/></w:rPr><w:t>Godwin</w:t></w:r><w:r w:rsidR="002C1399" w:rsidRPr="008C44EF"><w:rPr><w:color w:val="636363"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve"> Drive, Suite </w:t></w:r><w:r><w:rPr><w:color w:val="636363"/><w:szCs w:val="32"/></w:rPr><w:t>202</w:t></w:r><w:r w:rsidR="002C1399" w:rsidRPr="008C44EF"><w:rPr><w:color w:val="636363"/><w:szCs w:val="32"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:color w:val="636363"/><w:szCs w:val="32"/></w:rPr><w:t>Manassas</w:t></w:r><w:r w:rsidR="002C1399" w:rsidRPr="008C44EF"><w:rPr><w:color w:val="636363"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="gramStart"/><w:r w:rsidR="002C1399" w:rsidRPr="008C44EF"><w:rPr><w:color w:val="636363"/><w:szCs w:val="32"/></w:rPr><w:t>Virginia  2</w:t></w:r><w:r><w:rPr><w:color w:val="636363"/><w:szCs w:val="32"/></w:rPr><w:t>0110</w:t></w:r><w:proofErr w:type="gramEnd"/>'

# Locate paragraphs by their current text content (indices are stable since run/paragraph
# counts do not change -- we only ever add runs/proofErr markers inside existing paragraphs).
# Paragraph.Range.Text includes the trailing paragraph-mark (chr 13), so trim it first, and
# the embedded <w:br/> becomes chr 11 in Range.Text, so use a wildcard across it.
$idx1 = 0
$idx2 = 0
$idx3 = 0
$idx5 = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13)
    if ($t -eq "This is to certify that") { $idx1 = $i }
    if ($t -eq "Has successfully completed") { $idx2 = $i }
    if ($t -eq "Location LOCATION Date START_DATE END_DATE") { $idx3 = $i }
    if ($t -like "9817 Godwin Drive, Suite 202*Manassas, Virginia  20110") { $idx5 = $i }
}

if ($idx1 -eq 0) { throw "paragraph 1 (This is to certify that) not found" }
if ($idx2 -eq 0) { throw "paragraph 2 (Has successfully completed) not found" }
if ($idx3 -eq 0) { throw "paragraph 3 (Location LOCATION Date ...) not found" }
if ($idx5 -eq 0) { throw "paragraph 5 (9817 Godwin Drive ... 20110) not found" }

Replace-ParagraphRuns $idx1 $p1
Replace-ParagraphRuns $idx2 $p2
Replace-ParagraphRuns $idx3 $p3
Replace-ParagraphRuns $idx5 $p5

Write-Output "done"
